$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 6 (shifts "Polynomial Regression" ... "GB Regression" down by one)
$ws.Rows.Item(6).Insert()

# Copy the formatting from row 5 (Lasso Regression+normalization row) into the new row 6
$ws.Range("A5:C5").Copy()
$ws.Range("A6:C6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new row's data: Id=4, Model="Lasso Regression+normalization+ lag1", Accuracy=71.002143130205795
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "Lasso Regression+normalization+ lag1"
$ws.Range("C6").Value = 71.002143130205795

# Renumber the Id column for the rows that were pushed down (rows 7-11 now hold ids 5-9)
$ws.Range("A7").Value = 5
$ws.Range("A8").Value = 6
$ws.Range("A9").Value = 7
$ws.Range("A10").Value = 8
$ws.Range("A11").Value = 9

# Update the selection to match the saved workbook state
[void]$ws.Range("E5").Select()
